$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(3)
$tr = $shape.TextFrame.TextRange
$tr.Text = "x"
$tr.Text = "An image"
